$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking strings
# (e.g. "304.10", "1.01") are preserved verbatim instead of being coerced
# into floating point numbers (which would drop trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = "43.928.98"
$ws.Cells.Item(2, 5).Value = "  -0.49%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.231.80"
$ws.Cells.Item(3, 5).Value = "  -0.72%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "1.01"
$ws.Cells.Item(4, 5).Value = "  +0.33%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "304.10"
$ws.Cells.Item(5, 5).Value = "  -4.49%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "93.68"
$ws.Cells.Item(6, 5).Value = "  -7.39%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.87%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.25%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -4.49%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "34.58"
$ws.Cells.Item(10, 5).Value = "  -6.17%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "0.0805"
$ws.Cells.Item(11, 5).Value = "  -2.72%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "7.15"
$ws.Cells.Item(12, 5).Value = "  -4.81%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  -1.03%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "2.571.08"
$ws.Cells.Item(14, 5).Value = "  -0.74%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "2.233.30"
$ws.Cells.Item(15, 5).Value = "  -2.20%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "0.817"
$ws.Cells.Item(16, 5).Value = "  -3.71%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "13.46"
$ws.Cells.Item(17, 5).Value = "  -5.13%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "43.781.16"
$ws.Cells.Item(18, 5).Value = "  -0.53%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "0.0₃0957"
$ws.Cells.Item(19, 5).Value = "  -1.93%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "11.92"
$ws.Cells.Item(20, 5).Value = "  -11.43%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "6.24"
$ws.Cells.Item(21, 5).Value = "  -2.88%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "65.22"

# Row 23
$ws.Cells.Item(23, 4).Value = "236.03"
$ws.Cells.Item(23, 5).Value = "  +0.44%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "2.90"
$ws.Cells.Item(24, 5).Value = "  -6.54%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -5.55%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +0.51%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "9.83"
$ws.Cells.Item(27, 5).Value = "  -5.88%  "

# Row 28
$ws.Cells.Item(28, 2).Value = "Toncoin"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(28, 4).Value = "2.17"
$ws.Cells.Item(28, 5).Value = "  -2.12%  "

# Row 29
$ws.Cells.Item(29, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(29, 4).Value = "37.40"
$ws.Cells.Item(29, 5).Value = "  -0.30%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  -2.53%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "19.85"
$ws.Cells.Item(31, 5).Value = "  -1.15%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "150.93"
$ws.Cells.Item(32, 5).Value = "  -4.76%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "0.0794"
$ws.Cells.Item(33, 5).Value = "  -6.08%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "WEMIXToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(34, 4).Value = "2.58"
$ws.Cells.Item(34, 5).Value = "  -3.29%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "LidoDAOToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(35, 4).Value = "3.19"
$ws.Cells.Item(35, 5).Value = "  +0.71%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -3.00%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "0.119"
$ws.Cells.Item(37, 5).Value = "  +0.48%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -9.85%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "14.78"
$ws.Cells.Item(39, 5).Value = "  -8.09%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -8.40%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  -9.96%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "0.0294"
$ws.Cells.Item(42, 5).Value = "  -6.53%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  +0.37%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "1.724.88"
$ws.Cells.Item(44, 5).Value = "  -1.29%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "84.54"
$ws.Cells.Item(45, 5).Value = "  +2.78%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  -5.60%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "99.30"
$ws.Cells.Item(47, 5).Value = "  -3.34%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "4.90"
$ws.Cells.Item(48, 5).Value = "  -4.84%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "8.02"
$ws.Cells.Item(49, 5).Value = "  -3.07%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "EnergySwap"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(50, 4).Value = "14.32"
$ws.Cells.Item(50, 5).Value = "  +0.98%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "ordi"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Cells.Item(51, 4).Value = "67.97"
$ws.Cells.Item(51, 5).Value = "  -9.00%  "
